$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Menu'
$ws.Range('B2').Value = 15
$ws.Range('C2').Value = 5
$ws.Range('D2').Value = 20

$ws.Range('A3').Value = 'Tree'
$ws.Range('B3').Value = 10
$ws.Range('C3').Value = 0
$ws.Range('D3').Value = 10

$ws.Range('A4').Value = 'Nav'
$ws.Range('B4').Value = 8
$ws.Range('C4').Value = 1
$ws.Range('D4').Value = 9

$ws.Range('A5').Value = 'DataGrid'
$ws.Range('B5').Value = 3
$ws.Range('C5').Value = 6
$ws.Range('D5').Value = 9

$ws.Range('A6').Value = 'Popover'
$ws.Range('B6').Value = 6
$ws.Range('C6').Value = 2
$ws.Range('D6').Value = 8

$ws.Range('A7').Value = 'Dialog'
$ws.Range('B7').Value = 6
$ws.Range('C7').Value = 1
$ws.Range('D7').Value = 7

$ws.Range('A8').Value = 'Combobox'
$ws.Range('B8').Value = 3
$ws.Range('C8').Value = 3
$ws.Range('D8').Value = 6

$ws.Range('A9').Value = 'Table'
$ws.Range('B9').Value = 4
$ws.Range('C9').Value = 1
$ws.Range('D9').Value = 5

$ws.Range('A10').Value = 'Tooltip'
$ws.Range('B10').Value = 3
$ws.Range('C10').Value = 2
$ws.Range('D10').Value = 5

$ws.Range('A11').Value = 'Virtualizer'
$ws.Range('B11').Value = 4
$ws.Range('C11').Value = 0
$ws.Range('D11').Value = 4

$ws.Range('A12').Value = 'Dropdown'
$ws.Range('B12').Value = 2
$ws.Range('C12').Value = 2
$ws.Range('D12').Value = 4

$ws.Range('A13').Value = 'TagPicker'
$ws.Range('B13').Value = 3
$ws.Range('C13').Value = 1
$ws.Range('D13').Value = 4

$ws.Range('A14').Value = 'Toolbar'
$ws.Range('B14').Value = 3
$ws.Range('C14').Value = 1
$ws.Range('D14').Value = 4

$ws.Range('A15').Value = 'Skeleton'
$ws.Range('B15').Value = 2
$ws.Range('C15').Value = 1
$ws.Range('D15').Value = 3

$ws.Range('A16').Value = 'Calendar Compat'
$ws.Range('B16').Value = 3
$ws.Range('C16').Value = 0
$ws.Range('D16').Value = 3

$ws.Range('A17').Value = 'TeachingPopover'
$ws.Range('B17').Value = 3
$ws.Range('C17').Value = 0
$ws.Range('D17').Value = 3

$ws.Range('A18').Value = 'MessageBar'
$ws.Range('B18').Value = 3
$ws.Range('C18').Value = 0
$ws.Range('D18').Value = 3

$ws.Range('A19').Value = 'Switch'
$ws.Range('B19').Value = 2
$ws.Range('C19').Value = 0
$ws.Range('D19').Value = 2

$ws.Range('A20').Value = 'Toast'
$ws.Range('B20').Value = 1
$ws.Range('C20').Value = 1
$ws.Range('D20').Value = 2

$ws.Range('A21').Value = 'Drawer'
$ws.Range('B21').Value = 2
$ws.Range('C21').Value = 0
$ws.Range('D21').Value = 2

$ws.Range('A22').Value = 'Accordion'
$ws.Range('B22').Value = 2
$ws.Range('C22').Value = 0
$ws.Range('D22').Value = 2

$ws.Range('A23').Value = 'FluentProvider'
$ws.Range('B23').Value = 0
$ws.Range('C23').Value = 2
$ws.Range('D23').Value = 2

$ws.Range('A24').Value = 'DatePicker'
$ws.Range('B24').Value = 0
$ws.Range('C24').Value = 2
$ws.Range('D24').Value = 2

$ws.Range('A25').Value = 'Slider'
$ws.Range('B25').Value = 0
$ws.Range('C25').Value = 2
$ws.Range('D25').Value = 2

$ws.Range('A26').Value = 'Portal'
$ws.Range('B26').Value = 2
$ws.Range('C26').Value = 0
$ws.Range('D26').Value = 2

$ws.Range('A27').Value = 'Tabs'
$ws.Range('B27').Value = 2
$ws.Range('C27').Value = 0
$ws.Range('D27').Value = 2

$ws.Range('A28').Value = 'Popup'
$ws.Range('B28').Value = 1
$ws.Range('C28').Value = 0
$ws.Range('D28').Value = 1

$ws.Range('A29').Value = 'Input'
$ws.Range('B29').Value = 0
$ws.Range('C29').Value = 1
$ws.Range('D29').Value = 1

$ws.Range('A30').Value = 'Image'
$ws.Range('B30').Value = 0
$ws.Range('C30').Value = 1
$ws.Range('D30').Value = 1

$ws.Range('A31').Value = 'DatePickerCompat'
$ws.Range('B31').Value = 0
$ws.Range('C31').Value = 1
$ws.Range('D31').Value = 1

$ws.Range('A32').Value = 'Tag'
$ws.Range('B32').Value = 1
$ws.Range('C32').Value = 0
$ws.Range('D32').Value = 1

$ws.Range('A33').Value = 'MenuItem'
$ws.Range('B33').Value = 1
$ws.Range('C33').Value = 0
$ws.Range('D33').Value = 1

$ws.Range('A34').Value = 'AvatarGroup'
$ws.Range('B34').Value = 1
$ws.Range('C34').Value = 0
$ws.Range('D34').Value = 1

$ws.Range('A35').Value = 'Label'
$ws.Range('B35').Value = 1
$ws.Range('C35').Value = 0
$ws.Range('D35').Value = 1

$ws.Range('A36').Value = 'FocusTrapZone'
$ws.Range('B36').Value = 1
$ws.Range('C36').Value = 0
$ws.Range('D36').Value = 1

$ws.Range('A37').Value = 'Button'
$ws.Range('B37').Value = 0
$ws.Range('C37').Value = 1
$ws.Range('D37').Value = 1

$ws.Range('A38').Value = 'SwatchPicker'
$ws.Range('B38').Value = 0
$ws.Range('C38').Value = 1
$ws.Range('D38').Value = 1

$ws.Range('A39').Value = 'List'
$ws.Range('B39').Value = 1
$ws.Range('C39').Value = 0
$ws.Range('D39').Value = 1

$ws.Range('A40').Value = 'Avatar'
$ws.Range('B40').Value = 1
$ws.Range('C40').Value = 0
$ws.Range('D40').Value = 1

$ws.Range('A41').Value = 'Badge'
$ws.Range('B41').Value = 1
$ws.Range('C41').Value = 0
$ws.Range('D41').Value = 1

$ws.Range('A42').Value = 'Checkbox'
$ws.Range('B42').Value = 1
$ws.Range('C42').Value = 0
$ws.Range('D42').Value = 1

$ws.Range('A43').Value = 'SplitButton'
$ws.Range('B43').Value = 0
$ws.Range('C43').Value = 0
$ws.Range('D43').Value = 0

$ws.Range('A44').Value = 'InfoLabel'
$ws.Range('B44').Value = 0
$ws.Range('C44').Value = 0
$ws.Range('D44').Value = 0

$ws.Range('A45').Value = 'Rating'
$ws.Range('B45').Value = 0
$ws.Range('C45').Value = 0
$ws.Range('D45').Value = 0

$ws.Range('A46').Value = 'ColorPicker'
$ws.Range('B46').Value = 0
$ws.Range('C46').Value = 0
$ws.Range('D46').Value = 0

$ws.Range('A47').Value = 'SpinButton'
$ws.Range('B47').Value = 0
$ws.Range('C47').Value = 0
$ws.Range('D47').Value = 0

$ws.Range('A48').Value = 'Segment'
$ws.Range('B48').Value = 0
$ws.Range('C48').Value = 0
$ws.Range('D48').Value = 0

$ws.Range('A49').Value = 'Pickers'
$ws.Range('B49').Value = 0
$ws.Range('C49').Value = 0
$ws.Range('D49').Value = 0

$ws.Range('A50').Value = 'Keytip'
$ws.Range('B50').Value = 0
$ws.Range('C50').Value = 0
$ws.Range('D50').Value = 0
